# Workbook / sheet references
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 336, pushing all existing data
# (previously rows 336-408) down to rows 338-410.
$ws.Range("336:337").Insert()

# ---- New row 336 ----
$ws.Range("A336").Value = 5
$ws.Range("B336").Value = "Macroferia Regional de Talca"
$ws.Range("C336").Value = "Maule"
$ws.Range("D336").Value = 44511
$ws.Range("E336").Value = 7
$ws.Range("F336").Value = 100112004
$ws.Range("G336").Value = "Cebolla"
$ws.Range("H336").Value = "Sin especificar"
$ws.Range("I336").Value = "1a nueva(o)"
$ws.Range("J336").Value = 2000
$ws.Range("K336").Value = 3500
$ws.Range("L336").Value = 3500
$ws.Range("M336").Value = 3500
$ws.Range("N336").Value = "`$/malla 18 kilos"
$ws.Range("O336").Value = "Región de Arica y Parinacota"
$ws.Range("P336").Value = 194
$ws.Range("Q336").Value = 18
$ws.Range("R336").Value = "Hortaliza"

# ---- New row 337 ----
$ws.Range("A337").Value = 5
$ws.Range("B337").Value = "Macroferia Regional de Talca"
$ws.Range("C337").Value = "Maule"
$ws.Range("D337").Value = 44511
$ws.Range("E337").Value = 7
$ws.Range("F337").Value = 100112004
$ws.Range("G337").Value = "Cebolla"
$ws.Range("H337").Value = "Sin especificar"
$ws.Range("I337").Value = "1a nueva(o)"
$ws.Range("J337").Value = 50000
$ws.Range("K337").Value = 1200
$ws.Range("L337").Value = 1200
$ws.Range("M337").Value = 1200
$ws.Range("N337").Value = "`$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O337").Value = "Región de O'Higgins"
$ws.Range("P337").Value = 120
$ws.Range("Q337").Value = 10
$ws.Range("R337").Value = "Hortaliza"
